$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Monday (B6) hours for week commencing 43150 (row 6)
$ws.Range("B6").Value = 8.25

# Fill in Thursday (H6) hours for that week, which was previously blank
$ws.Range("H6").Value = 6.25

# Update the currently selected cell/range to E6
$ws.Range("E6").Select()

$wb.Save()
